$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').NumberFormat = '@'
$ws.Range('A2').Value = '281474990868838-1749474760673'
$ws.Range('A2').Style = 'Normal'

$ws.Range('B2').NumberFormat = '@'
$ws.Range('B2').Value = 'Mobile Usage'
$ws.Range('B2').Style = 'Normal'

$ws.Range('C2').NumberFormat = '@'
$ws.Range('C2').Value = '2025-06-09T07:12:40.673'
$ws.Range('C2').Style = 'Normal'

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '281474990868838'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '135'
$ws.Range('E2').Style = 'Normal'

$ws.Range('F2').NumberFormat = '@'
$ws.Range('F2').Value = '51834045'
$ws.Range('F2').Style = 'Normal'

$ws.Range('G2').NumberFormat = '@'
$ws.Range('G2').Value = 'FRANCISCO JAVIER MÁRQUEZ PÉREZ'
$ws.Range('G2').Style = 'Normal'

$ws.Range('H2').Value = 19.88336851
$ws.Range('H2').Style = 'Normal'

$ws.Range('I2').Value = -101.43772428
$ws.Range('I2').Style = 'Normal'

$ws.Range('J2').Value = 0
$ws.Range('J2').Style = 'Normal'

$ws.Range('K2').NumberFormat = '@'
$ws.Range('K2').Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990868838/1749474758173/gFRHQATQ91-camera-video-segment-driver-1749474760673.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSAHEZKD52%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T150652Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQCsTJzRA78IKu7CIWoXR0acB5P%2BxvHg0Tv%2FTO42eviyDAIhAN7E4GFk%2FNOMftxQ77jH2R1Mo0Sx%2FOSePkSR7D9A46hbKuYDCL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0Igx%2FXAb%2FzrRUm0ummDwqugOLEmrJQ6HrKTcviHpEVbF339krKtEh643Pw5m2eg0Glj9cdasveFURkRcutgfRr16CPDTn%2BLNXg3ry%2BlIijo020ajuK44LiRx%2FcOaps3BtPF4vk0vbh0ahNfApWJv3jdGWmDpjsKsfaNmGH3hDdQzlBbQgnkCGnSOXWiW6oRV5sTQPf0A7t1%2FsmXpq%2ByEYUdivZeONAWVL%2FL3rFPO6HVyv4wQXchm2uLb4ZLbgz0Pem9mnTtD%2BQUS3h4eNXCFK2eFZkwH75ASTCLweLGXbDPFRkpvtnso178u5MQkp%2B0W4vKSwASnCJ16SOkxE4RC3lCz5wLvuOHTNNPoY5NEE6oLHFbCQEIt0oq%2F4Hjme2oM0ojJ0vw1H7hJs%2B45Fn8Vs0fmZqI2JLCSShQ1nfSx0MtQ9b5Bz6g6u5Z7opdYooM02J5XllVdEi7KA0FjpPAI6td%2B1n3Oaz1j48zCe%2FwAFGIri7VLtwkabHlEzKNiH5yazCCRdxZwU%2FcKn2JlammGoUR54yJh7WqM6Y8J3izjc%2Bqzr8XyHEHr0n%2FSC1PvL0qIG0XZlILiG2DNx1BLKkHvi5oE5vsai9d4WuwcCMPDtoMIGOqQBcMbafWad9t%2B3%2FygezrA1hYAlcvwkQvNu5ReY1zxtS7TbZgAj9UV1NsA00aNwFLps27wWjfQqHJvseZvb5OglMYGG1zU6BDj5F6sGGKFqK%2FsAiYVmxPgQoQu8v3SDSGgJbTURQddVmVNWWZv%2BEupyPV5735F67LWzL6Wrs1yIBe6T5QoNR6vkAn40d%2Fou3zSXwQaH90NZz9mAg9QaK%2B8MtEyuNN8%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2010%20Jun%202025%2023%3A06%3A52%20GMT&X-Amz-Signature=d85c79b9dd19eab6e80e6dba09ffa144479948c68b7d482b5a3e4dce15b438a9'
$ws.Range('K2').Style = 'Normal'

$ws.Range('L2').NumberFormat = '@'
$ws.Range('L2').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990868838/1749474758173/7Z3ova6Wtr-camera-video-segment-1749474760673.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSAHEZKD52%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T150652Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQCsTJzRA78IKu7CIWoXR0acB5P%2BxvHg0Tv%2FTO42eviyDAIhAN7E4GFk%2FNOMftxQ77jH2R1Mo0Sx%2FOSePkSR7D9A46hbKuYDCL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0Igx%2FXAb%2FzrRUm0ummDwqugOLEmrJQ6HrKTcviHpEVbF339krKtEh643Pw5m2eg0Glj9cdasveFURkRcutgfRr16CPDTn%2BLNXg3ry%2BlIijo020ajuK44LiRx%2FcOaps3BtPF4vk0vbh0ahNfApWJv3jdGWmDpjsKsfaNmGH3hDdQzlBbQgnkCGnSOXWiW6oRV5sTQPf0A7t1%2FsmXpq%2ByEYUdivZeONAWVL%2FL3rFPO6HVyv4wQXchm2uLb4ZLbgz0Pem9mnTtD%2BQUS3h4eNXCFK2eFZkwH75ASTCLweLGXbDPFRkpvtnso178u5MQkp%2B0W4vKSwASnCJ16SOkxE4RC3lCz5wLvuOHTNNPoY5NEE6oLHFbCQEIt0oq%2F4Hjme2oM0ojJ0vw1H7hJs%2B45Fn8Vs0fmZqI2JLCSShQ1nfSx0MtQ9b5Bz6g6u5Z7opdYooM02J5XllVdEi7KA0FjpPAI6td%2B1n3Oaz1j48zCe%2FwAFGIri7VLtwkabHlEzKNiH5yazCCRdxZwU%2FcKn2JlammGoUR54yJh7WqM6Y8J3izjc%2Bqzr8XyHEHr0n%2FSC1PvL0qIG0XZlILiG2DNx1BLKkHvi5oE5vsai9d4WuwcCMPDtoMIGOqQBcMbafWad9t%2B3%2FygezrA1hYAlcvwkQvNu5ReY1zxtS7TbZgAj9UV1NsA00aNwFLps27wWjfQqHJvseZvb5OglMYGG1zU6BDj5F6sGGKFqK%2FsAiYVmxPgQoQu8v3SDSGgJbTURQddVmVNWWZv%2BEupyPV5735F67LWzL6Wrs1yIBe6T5QoNR6vkAn40d%2Fou3zSXwQaH90NZz9mAg9QaK%2B8MtEyuNN8%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2010%20Jun%202025%2023%3A06%3A52%20GMT&X-Amz-Signature=716da3e986533395ea365e1a82f69a54d3f3ca121eb41cd8c98c28831e1f2d5a'
$ws.Range('L2').Style = 'Normal'

$ws.Range('A3').NumberFormat = '@'
$ws.Range('A3').Value = '281474990868838-1749474691113'
$ws.Range('A3').Style = 'Normal'

$ws.Range('B3').NumberFormat = '@'
$ws.Range('B3').Value = 'Mobile Usage'
$ws.Range('B3').Style = 'Normal'

$ws.Range('C3').NumberFormat = '@'
$ws.Range('C3').Value = '2025-06-09T07:11:31.113'
$ws.Range('C3').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '281474990868838'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '135'
$ws.Range('E3').Style = 'Normal'

$ws.Range('F3').NumberFormat = '@'
$ws.Range('F3').Value = '51834045'
$ws.Range('F3').Style = 'Normal'

$ws.Range('G3').NumberFormat = '@'
$ws.Range('G3').Value = 'FRANCISCO JAVIER MÁRQUEZ PÉREZ'
$ws.Range('G3').Style = 'Normal'

$ws.Range('H3').Value = 19.88523483
$ws.Range('H3').Style = 'Normal'

$ws.Range('I3').Value = -101.4518704
$ws.Range('I3').Style = 'Normal'

$ws.Range('J3').Value = 0
$ws.Range('J3').Style = 'Normal'

$ws.Range('K3').NumberFormat = '@'
$ws.Range('K3').Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990868838/1749474688613/JLWgs3viB6-camera-video-segment-driver-1749474691113.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSAHEZKD52%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T150652Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQCsTJzRA78IKu7CIWoXR0acB5P%2BxvHg0Tv%2FTO42eviyDAIhAN7E4GFk%2FNOMftxQ77jH2R1Mo0Sx%2FOSePkSR7D9A46hbKuYDCL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0Igx%2FXAb%2FzrRUm0ummDwqugOLEmrJQ6HrKTcviHpEVbF339krKtEh643Pw5m2eg0Glj9cdasveFURkRcutgfRr16CPDTn%2BLNXg3ry%2BlIijo020ajuK44LiRx%2FcOaps3BtPF4vk0vbh0ahNfApWJv3jdGWmDpjsKsfaNmGH3hDdQzlBbQgnkCGnSOXWiW6oRV5sTQPf0A7t1%2FsmXpq%2ByEYUdivZeONAWVL%2FL3rFPO6HVyv4wQXchm2uLb4ZLbgz0Pem9mnTtD%2BQUS3h4eNXCFK2eFZkwH75ASTCLweLGXbDPFRkpvtnso178u5MQkp%2B0W4vKSwASnCJ16SOkxE4RC3lCz5wLvuOHTNNPoY5NEE6oLHFbCQEIt0oq%2F4Hjme2oM0ojJ0vw1H7hJs%2B45Fn8Vs0fmZqI2JLCSShQ1nfSx0MtQ9b5Bz6g6u5Z7opdYooM02J5XllVdEi7KA0FjpPAI6td%2B1n3Oaz1j48zCe%2FwAFGIri7VLtwkabHlEzKNiH5yazCCRdxZwU%2FcKn2JlammGoUR54yJh7WqM6Y8J3izjc%2Bqzr8XyHEHr0n%2FSC1PvL0qIG0XZlILiG2DNx1BLKkHvi5oE5vsai9d4WuwcCMPDtoMIGOqQBcMbafWad9t%2B3%2FygezrA1hYAlcvwkQvNu5ReY1zxtS7TbZgAj9UV1NsA00aNwFLps27wWjfQqHJvseZvb5OglMYGG1zU6BDj5F6sGGKFqK%2FsAiYVmxPgQoQu8v3SDSGgJbTURQddVmVNWWZv%2BEupyPV5735F67LWzL6Wrs1yIBe6T5QoNR6vkAn40d%2Fou3zSXwQaH90NZz9mAg9QaK%2B8MtEyuNN8%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2010%20Jun%202025%2023%3A06%3A52%20GMT&X-Amz-Signature=e65328c5ff134478a0608c02c0248ad90deb8997475bbe4b15414f0c9d98483c'
$ws.Range('K3').Style = 'Normal'

$ws.Range('L3').NumberFormat = '@'
$ws.Range('L3').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990868838/1749474688613/LOHXs2uAgt-camera-video-segment-1749474691113.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSAHEZKD52%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T150652Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQCsTJzRA78IKu7CIWoXR0acB5P%2BxvHg0Tv%2FTO42eviyDAIhAN7E4GFk%2FNOMftxQ77jH2R1Mo0Sx%2FOSePkSR7D9A46hbKuYDCL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0Igx%2FXAb%2FzrRUm0ummDwqugOLEmrJQ6HrKTcviHpEVbF339krKtEh643Pw5m2eg0Glj9cdasveFURkRcutgfRr16CPDTn%2BLNXg3ry%2BlIijo020ajuK44LiRx%2FcOaps3BtPF4vk0vbh0ahNfApWJv3jdGWmDpjsKsfaNmGH3hDdQzlBbQgnkCGnSOXWiW6oRV5sTQPf0A7t1%2FsmXpq%2ByEYUdivZeONAWVL%2FL3rFPO6HVyv4wQXchm2uLb4ZLbgz0Pem9mnTtD%2BQUS3h4eNXCFK2eFZkwH75ASTCLweLGXbDPFRkpvtnso178u5MQkp%2B0W4vKSwASnCJ16SOkxE4RC3lCz5wLvuOHTNNPoY5NEE6oLHFbCQEIt0oq%2F4Hjme2oM0ojJ0vw1H7hJs%2B45Fn8Vs0fmZqI2JLCSShQ1nfSx0MtQ9b5Bz6g6u5Z7opdYooM02J5XllVdEi7KA0FjpPAI6td%2B1n3Oaz1j48zCe%2FwAFGIri7VLtwkabHlEzKNiH5yazCCRdxZwU%2FcKn2JlammGoUR54yJh7WqM6Y8J3izjc%2Bqzr8XyHEHr0n%2FSC1PvL0qIG0XZlILiG2DNx1BLKkHvi5oE5vsai9d4WuwcCMPDtoMIGOqQBcMbafWad9t%2B3%2FygezrA1hYAlcvwkQvNu5ReY1zxtS7TbZgAj9UV1NsA00aNwFLps27wWjfQqHJvseZvb5OglMYGG1zU6BDj5F6sGGKFqK%2FsAiYVmxPgQoQu8v3SDSGgJbTURQddVmVNWWZv%2BEupyPV5735F67LWzL6Wrs1yIBe6T5QoNR6vkAn40d%2Fou3zSXwQaH90NZz9mAg9QaK%2B8MtEyuNN8%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2010%20Jun%202025%2023%3A06%3A52%20GMT&X-Amz-Signature=4327639419774862b7a8f89b9400a48122639339944972f73af3a01fe84467c4'
$ws.Range('L3').Style = 'Normal'
